$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 12: Helsinki University Hospital ---
$ws.Range("B12").Value = 48
$ws.Range("D12").Value = 4.2
$ws.Range("F12").Value = 14

# --- Row 21: Lund University ---
$ws.Range("E21").Value = 4.399999999999999

# --- Rows 27-55: re-insert "Örebro University" / "Örebro University Hospital"
# right after "Odense University Hospital" (row 26), pushing the
# intervening institutions down by two rows. Rebuild the block in the
# new order, data only (no formatting changes needed - these rows carry
# the default style).

$data = @(
    @("Örebro University", 17, 0, 0, 0, 18.4),
    @("Örebro University Hospital", 1, 0, 0, 0, 94.89999999999999),
    @("Oslo University Hospital", 102, 1, 1, 0.1, 5.3),
    @("Oulu University Hospital", 10, 0, 0, 0, 27.8),
    @("Sahlgrenska University Hospital", 42, 0, 0, 0, 8.4),
    @("Skane University Hospital", 23, 0, 0, 0, 14.3),
    @("St. Olav’s University Hospital", 24, 0, 0, 0, 13.8),
    @("Steno Diabetes Center Copenhagen", 13, 0, 0, 0, 22.8),
    @("Stockholm South General Hospital", 3, 0, 0, 0, 56.10000000000001),
    @("Tampere University Hospital", 24, 0, 0, 0, 13.8),
    @("The National University Hospital of Iceland", 5, 0, 0, 0, 43.4),
    @("Turku University Hospital", 50, 0, 0, 0, 7.1),
    @("UiT The Arctic University of Norway", 14, 0, 0, 0, 21.5),
    @("Umeå University", 42, 0, 0, 0, 8.4),
    @("University Hospital of North Norway", 17, 0, 0, 0, 18.4),
    @("University Hospital of Umeå", 2, 1, 50, 2.6, 97.39999999999999),
    @("University of Bergen", 31, 0, 0, 0, 11),
    @("University of Copenhagen", 99, 1, 1, 0.1, 5.5),
    @("University of Eastern Finland", 12, 0, 0, 0, 24.2),
    @("University of Helsinki", 21, 0, 0, 0, 15.5),
    @("University of Iceland", 5, 0, 0, 0, 43.4),
    @("University of Oslo", 23, 0, 0, 0, 14.3),
    @("University of Oulu", 25, 1, 4, 0.2, 19.5),
    @("University of Southern Denmark", 42, 0, 0, 0, 8.4),
    @("University of Tampere", 9, 1, 11.1, 0.6, 43.5),
    @("University of Turku", 20, 0, 0, 0, 16.1),
    @("Uppsala Academic Hospital", 9, 0, 0, 0, 29.9),
    @("Uppsala University", 51, 1, 2, 0.1, 10.3),
    @("Zealand University Hospital", 28, 0, 0, 0, 12.1)
)

$r = 27
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
    $ws.Cells.Item($r, 6).Value = $row[5]
    $r++
}

# --- Row 56: Total trials count ---
$ws.Range("B56").Value = 2112
